$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update existing sheet "o_10": refresh prompt/solution/llm_response text for the new question ---

$promptText = @"
 Given is the adjacency matrix for a weighted undirected graph containing 17 nodes labelled A to Q. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: What is the least cost path to travel first from node A to node I, and then from node node I to node J?
   A B C D E F G H I J
 A 0 1 0 0 2 0 3 0 0 5
 B 1 0 3 4 0 0 0 0 0 0
 C 0 3 0 0 0 0 0 0 0 0
 D 0 4 0 0 0 0 0 0 0 0
 E 2 0 0 0 0 2 0 0 0 0
 F 0 0 0 0 2 0 0 0 0 0
 G 3 0 0 0 0 0 0 1 0 0
 H 0 0 0 0 0 0 1 0 3 0
 I 0 0 0 0 0 0 0 3 0 4
 J 5 0 0 0 0 0 0 0 4 0

Solution: Path from A to I: A -> G -> H -> I
Path from I to J: I -> J
        

Example 2: What is the least cost path to travel first from node A to node B, and then from node node B to node K?
   A B C D E F G H I J K
 A 0 1 3 0 0 0 0 0 1 0 4
 B 1 0 0 0 0 0 0 0 0 0 0
 C 3 0 0 4 0 0 3 2 0 0 0
 D 0 0 4 0 5 0 0 0 0 0 0
 E 0 0 0 5 0 3 0 0 0 0 0
 F 0 0 0 0 3 0 0 0 0 0 0
 G 0 0 3 0 0 0 0 0 0 0 0
 H 0 0 2 0 0 0 0 0 0 0 0
 I 1 0 0 0 0 0 0 0 0 5 0
 J 0 0 0 0 0 0 0 0 5 0 3
 K 4 0 0 0 0 0 0 0 0 3 0

Solution: Path from A to B: A -> B
Path from B to K: B -> A -> K
        

Example 3: What is the least cost path to travel first from node A to node K, and then from node node K to node S?
   A B C D E F G H I J K L M N O P Q R S
 A 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2
 B 3 0 5 3 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 4 0 0 0 5 2 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 2 0 0 4 0 0 2 5 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 4 0 2 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 2 0 1 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 5 0 0 0 0 0 5 0 4 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 5 0 4 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 2 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 2
 S 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0

Solution: Path from A to K: A -> B -> E -> G -> K
Path from K to S: K -> G -> E -> B -> A -> S
        
 Given these examples, answer the following quesiton.

What is the least cost path to travel first from node 0 to node I, and then from node node I to node 16?
   A B C D E F G H I J K L M N O P Q
 A 0 4 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 B 4 0 2 5 4 0 0 0 0 0 0 0 5 0 0 0 0
 C 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 4 0 0 0 2 3 0 0 0 0 0 0 0 0 0 5
 F 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 3 0 0 2 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 2 0 4 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 4 0 3 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 3 0 0 5 0 0 0 0 0
 K 1 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 5 4 0 0 0 0 0 0
 M 0 5 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 3 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 2
 Q 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 2 0
    
"@

$solutionText = @"
Path from A to I: A -> K -> L -> J -> I
Path from I to Q: I -> H -> G -> E -> Q
"@

$llmResponseText = @"
The least cost path to travel first from node 0 to node I, and then from node I to node 16 is:
0 -> B -> E -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> 16
"@

$ws1.Range("A2").Value = $promptText
$ws1.Range("B2").Value = $solutionText
$ws1.Range("C2").Value = $llmResponseText
$ws1.Rows.Item(2).AutoFit()

# --- Add the new "evaluator_partial_correctness" column (E), copying the header style from D1 ---

$ws1.Range("D1").Copy($ws1.Range("E1"))
$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("E2").Value = "N/A"

# --- Add new sheets "o_20" and "o_20_jumbled" after "o_10", each with the same 5-column header row ---

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "o_20"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "o_20_jumbled"

$ws1.Range("A1:E1").Copy($ws2.Range("A1"))
$ws1.Range("A1:E1").Copy($ws3.Range("A1"))

$ws1.Select()

Write-Host "edit complete"
